$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 84; existing rows 84-118 shift down to 85-119.
$ws.Rows.Item(84).Insert()

# Populate the newly inserted row 84 with the new record's data.
$ws.Cells.Item(84, 1).Value = 9
$ws.Cells.Item(84, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(84, 3).Value = "Metropolitana"
$ws.Cells.Item(84, 4).Value = 44510
$ws.Cells.Item(84, 5).Value = 13
$ws.Cells.Item(84, 6).Value = "Fruta"
$ws.Cells.Item(84, 7).Value = 100101
$ws.Cells.Item(84, 8).Value = "Berries"
$ws.Cells.Item(84, 9).Value = 100101001
$ws.Cells.Item(84, 10).Value = "Arándano (blue)"
$ws.Cells.Item(84, 11).Value = "Sin especificar"
$ws.Cells.Item(84, 12).Value = "Primera"
$ws.Cells.Item(84, 13).Value = 300
$ws.Cells.Item(84, 14).Value = 8000
$ws.Cells.Item(84, 15).Value = 8000
$ws.Cells.Item(84, 16).Value = 8000
$ws.Cells.Item(84, 17).Value = "`$/bandeja 2 kilos"
$ws.Cells.Item(84, 18).Value = "Región Metropolitana"
$ws.Cells.Item(84, 19).Value = 4000
$ws.Cells.Item(84, 20).Value = 2
